$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 585, shifting existing rows 585:626 down to 586:627
$ws.Rows.Item(585).Insert()

$ws.Cells.Item(585, 1).Value = "2026/01/07"
$ws.Cells.Item(585, 2).Value = "水"
$ws.Cells.Item(585, 3).Value = 14
$ws.Cells.Item(585, 4).Value = 201
